# Apply quant engine data update: add Jan_2026 column, shift Dec_2025/Nov_2025,
# re-rank holdings by new Jan_2026 weight, add newly-entered holdings, append exited holdings at bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update month column headers: D/E/F shift from Dec/Nov/Oct_2025 to Jan_2026/Dec_2025/Nov_2025
$ws.Cells.Item(1, 4).Value = "Jan_2026"
$ws.Cells.Item(1, 5).Value = "Dec_2025"
$ws.Cells.Item(1, 6).Value = "Nov_2025"

# Full replacement holdings table (sorted descending by Jan_2026 weight, as produced by quant engine)
$rows = @(
    ,@("INE040A01034", "HDFC Bank Limited", 9.187924, 0, 0, 9.187924, 9.187924)
    ,@("INE202B01038", "Piramal Finance Ltd", 8.361349, 10.172964, 9.2446, -1.811615, -0.8832509999999996)
    ,@("INE018A01030", "Larsen & Toubro Limited", 7.955206, 7.788444, 7.483086, 0.1667620000000003, 0.4721200000000003)
    ,@("INE423A01024", "Adani Enterprises Limited", 7.508014, 7.84677, 7.812093, -0.3387560000000001, -0.3040789999999998)
    ,@("INE090A01021", "ICICI Bank Limited", 6.406086, 0, 0, 6.406086, 6.406086)
    ,@("INE795G01014", "HDFC Life Insurance Co Ltd", 5.947946, 5.752245, 5.652463, 0.1957009999999997, 0.2954829999999999)
    ,@("INE364U01010", "Adani Green Energy Limited", 5.256688, 5.902581, 3.887941, -0.645893, 1.368746999999999)
    ,@("INE406A01037", "Aurobindo Pharma Limited", 3.961293, 3.658272, 3.657131, 0.3030209999999998, 0.3041619999999998)
    ,@("INE917I01010", "Bajaj Auto Limited", 3.561141, 3.268366, 3.060061, 0.2927750000000002, 0.50108)
    ,@("INE237A01036", "Kotak Mahindra Bank Limited", 3.097165, 0, 0, 3.097165, 3.097165)
    ,@("INE814H01029", "Adani Power Limited", 2.788864, 4.034234, 4.011704, -1.24537, -1.22284)
    ,@("INE726G01019", "ICICI Prudential Life Insurance Co Ltd", 2.351255, 0, 0, 2.351255, 2.351255)
    ,@("INE200M01039", "Varun Beverages Limited", 1.93914, 0, 0, 1.93914, 1.93914)
    ,@("INE259A01022", "Colgate-Palmolive (India) Ltd", 1.572049, 0, 0, 1.572049, 1.572049)
    ,@("INE931S01010", "Adani Energy Solutions Limited", 0.423678, 0, 0, 0.423678, 0.423678)
    ,@("INE271C01023", "DLF Limited", 0, 5.824721, 5.911192, -5.824721, -5.911192)
    ,@("INE296A01032", "Bajaj Finance Limited", 0, 0, 6.334069, 0, -6.334069)
    ,@("INE237A01028", "Kotak Mahindra Bank Limited", 0, 3.150556, 2.931532, -3.150556, -2.931532)
    ,@("INE423A20016", "Adani Enterprises Limited Rights", 0, 0.18188, 0, 0, -0.18188)
    ,@("INE437A01024", "Apollo Hospitals Enterprise Ltd", 0, 1.997156, 0, -1.997156, 0)
    ,@("INE758E01017", "Jio Financial Services Limited", 0, 8.168625, 8.175543, -8.168625, -8.175543)
    ,@("INE775A01035", "Samvardhana Motherson International Ltd", 0, 6.446333, 3.129867, -6.446333, -3.129867)
    ,@("INE860A01027", "HCL Technologies Limited", 0, 1.618943, 0, -1.618943, 0)
)

$fundName = "quant Quantamental Fund"
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $fundName
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

"Wrote " + ($r - 2) + " data rows; last row = " + ($r - 1)
